# ---------------------------------------------------------------------------
# "supercambio" edit: flip the sign convention of the local/global stiffness
# sub-matrices (K1 K_altura, Barra 1 K1, Global) on the off-diagonal-style
# blocks, and refresh the dependent results (F, Desplazamientos, Reacciones)
# with the recomputed values.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Flip-MatrixSheet($ws, $n, $axialSet) {
    for ($r = 1; $r -le $n; $r++) {
        $rInSet = $axialSet -contains $r
        for ($c = 1; $c -le $n; $c++) {
            $cInSet = $axialSet -contains $c
            if ($rInSet -ne $cInSet) {
                $cell = $ws.Cells.Item($r, $c)
                $v = $cell.Value2
                $cell.Value2 = -$v
            }
        }
    }
}

# --- sheet "K1 (K_altura)": 12x12 local-axis stiffness matrix -------------
$wsAltura = $wb.Worksheets.Item("K1 (K_altura)")
Flip-MatrixSheet $wsAltura 12 @(1,4,7,10)

# --- sheet "Barra 1 (K1)": 12x12 stiffness matrix (same layout) -----------
$wsBarra1K1 = $wb.Worksheets.Item("Barra 1 (K1)")
Flip-MatrixSheet $wsBarra1K1 12 @(1,4,7,10)

# --- sheet "Global": 24x24 assembled stiffness matrix ----------------------
$wsGlobal = $wb.Worksheets.Item("Global")
Flip-MatrixSheet $wsGlobal 24 @(1,4,7,10,13,16,19,22)

# --- sheet "F": recomputed load vector --------------------------------------
$wsF = $wb.Worksheets.Item("F")
$fValues = @(
    "2.748314368632023e-16",
    "1.75973421993209",
    "-1.999697977195557",
    "2.399637572634668",
    "-2.077980688793229e-16",
    "1.469354236220385e-16",
    "-3.5",
    "3.239510723056801",
    "-14.51664747064227",
    "-3.599456358952002",
    "20.30570230503482",
    "-1.56139138049941e-15",
    "3.5",
    "-0.5199214740708445",
    "2.51725151625115",
    "0.7998791908782228",
    "-5.756475521456251",
    "2.633064592418956e-16",
    "-6.327821012781338e-16",
    "-4.479323468918047",
    "-3.999395954391114",
    "-3.199516763512891",
    "2.770640918390973e-16",
    "1.959138981627181e-16"
)
for ($i = 0; $i -lt $fValues.Length; $i++) {
    $wsF.Cells.Item($i + 1, 1).Value2 = [double]$fValues[$i]
}

# --- sheet "Desplazamientos": recomputed displacement vector --------------
$wsDesp = $wb.Worksheets.Item("Desplazamientos")
$despValues = @(
    "0.01380814379189671",
    "0.02167331059878497",
    "-0.2126831458241319",
    "0.00537055645510233",
    "-0.006031248225553358",
    "-0.002183288066494276",
    "0.05124197414119706",
    "-0.0004481385117980237",
    "-0.07301675785184697",
    "-8.457624555846062e-05",
    "-0.01772117031976563",
    "-0.002183288066494276"
)
for ($i = 0; $i -lt $despValues.Length; $i++) {
    $wsDesp.Cells.Item($i + 7, 1).Value2 = [double]$despValues[$i]
}

# --- sheet "Reacciones": recomputed reaction matrix (3x12) -----------------
$wsReac = $wb.Worksheets.Item("Reacciones")
$reacValues = @(
    @("2.713889562664689","-4.985137408818432","10.9323901018091","-14.8390749444679","20.06639319350182","0.07053767085229616","-2.71388956266469","-0.01410753417045951","-5.933145158820207","-0.08812221364647188","-33.63584100682528","-0.07053767085229709"),
    @("2.713889562664693","0.01410753417046012","5.933145158820203","0.08812221364647442","33.63584100682527","0.07053767085229688","-2.713889562664693","-0.01410753417045963","2.066854841179799","-0.08812221364647442","7.654485669883327","0.07053767085229759"),
    @("2.713889562664697","0.01410753417045951","-2.066854841179794","0.08812221364647432","-7.654485669883286","-0.07053767085229706","-2.713889562664697","4.985137408818433","7.066099784168687","4.84058505849012","21.22393348320678","0.07053767085229677")
)
for ($r = 0; $r -lt $reacValues.Length; $r++) {
    for ($c = 0; $c -lt $reacValues[$r].Length; $c++) {
        $wsReac.Cells.Item($r + 1, $c + 1).Value2 = [double]$reacValues[$r][$c]
    }
}
